$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: the existing "Segunda" entry is superseded by a new weekly
#     "Primera" observation (date moves forward to 2022-12-28) ---
$ws.Range("D3").Value = 44923
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 7500
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 7625
$ws.Range("S3").Value = 7625

# --- Row 4 (new): duplicate weekly "Primera" observation, same date/price ---
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44923
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7625
$ws.Range("Q4").Value = "$/bandeja 2,5 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 7625
$ws.Range("T4").Value = 1

# --- Row 5 (new): the original "Segunda" observation, re-appended ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44881
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101001
$ws.Range("J5").Value = "Arándano (blue)"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 11250
$ws.Range("O5").Value = 11250
$ws.Range("P5").Value = 11250
$ws.Range("Q5").Value = "$/bandeja 2,5 kilos"
$ws.Range("R5").Value = "Región de Coquimbo"
$ws.Range("S5").Value = 11250
$ws.Range("T5").Value = 1

# Match the "Fecha" column's date number format used by row 3 on the two
# freshly-created rows so D4/D5 render as dates like D2/D3.
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("D5").NumberFormat = $ws.Range("D3").NumberFormat
